# Converts a target EMU value into the point value that, once this
# runtime stores it as a 32-bit float (as PowerPoint's Shape geometry
# does) and converts back to EMU, reproduces that exact EMU value.
# (1 pt = 12700 EMU; plain division/round-tripping through float32 is
# lossy, so we search nearby doubles instead of trusting one division.)
function EmuToPt($targetEmu) {
    $emuPerPt = 12700.0
    $base = $targetEmu / $emuPerPt
    for ($i = -2000; $i -le 2000; $i++) {
        $cand = $base + ($i * 0.0000001)
        $f32 = [float]$cand
        $emu = [double]$f32 * $emuPerPt
        $floored = [math]::Floor($emu + 0.0000001)
        if ($floored -eq $targetEmu) {
            return $cand
        }
    }
    return $base
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The target shape must receive id=6 / auto-name "TextBox 5". The object
# model assigns shape ids/names from a monotonically increasing per-session
# counter that skips ids already used in the slide (it never reuses an id,
# even across deletes) - this mirrors the same "add, delete, add" pattern
# recorded for this slide (shape id 4 was added & removed, shape id 6 was
# added & kept). Reproduce that so the new textbox lands on id 6.
$discard = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$discard.Delete()

$left   = EmuToPt 317358
$top    = EmuToPt 4918770
$width  = EmuToPt 5022574
$height = EmuToPt 646331

$shp = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$shp.Name = "TextBox 5"

$tr = $shp.TextFrame.TextRange
$tr.Text = "https://github.com/ksu-is/brand-blog"

# Add a trailing empty paragraph (matches the authored slide, where the
# user pressed Enter after the pasted link).
[void]$tr.InsertAfter([char]13)

# Turn the URL text into a real hyperlink (first paragraph only).
$urlRange = $tr.Paragraphs(1, 1)
$urlRange.ActionSettings.Item(1).Hyperlink.Address = "https://github.com/ksu-is/brand-blog"

# Match the authored body formatting: word-wrap on, shape auto-fits to text.
$shp.TextFrame.WordWrap = -1
$shp.TextFrame.AutoSize = 1

# Re-assert the exact target geometry, since AutoSize recalculates it from
# the (possibly slightly different) rendered text metrics.
$shp.Left = $left
$shp.Top = $top
$shp.Width = $width
$shp.Height = $height
